# Update "想去人数" (want-to-go count) values in F column across sheets,
# reflecting a refreshed data scrape (gh-pages output regenerated).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet updates
$wsExhibit.Range("F8").Value  = 12966
$wsExhibit.Range("F9").Value  = 2247
$wsExhibit.Range("F11").Value = 290
$wsExhibit.Range("F12").Value = 53208
$wsExhibit.Range("F14").Value = 285
$wsExhibit.Range("F16").Value = 845
$wsExhibit.Range("F20").Value = 833
$wsExhibit.Range("F21").Value = 4978
$wsExhibit.Range("F22").Value = 1222
$wsExhibit.Range("F28").Value = 1170
$wsExhibit.Range("F30").Value = 17
$wsExhibit.Range("F37").Value = 4633
$wsExhibit.Range("F39").Value = 4693
$wsExhibit.Range("F40").Value = 5642
$wsExhibit.Range("F47").Value = 62
$wsExhibit.Range("F48").Value = 4148
$wsExhibit.Range("F49").Value = 163

# 演出 (Performance) sheet updates
$wsShow.Range("F12").Value = 1090

# 全部类型 (All types) sheet updates
$wsAll.Range("F7").Value  = 12966
$wsAll.Range("F8").Value  = 12966
$wsAll.Range("F9").Value  = 2247
$wsAll.Range("F10").Value = 290
$wsAll.Range("F13").Value = 845
$wsAll.Range("F17").Value = 833
$wsAll.Range("F19").Value = 4978
$wsAll.Range("F20").Value = 1222
$wsAll.Range("F26").Value = 1170
$wsAll.Range("F34").Value = 4633
$wsAll.Range("F36").Value = 4693
$wsAll.Range("F37").Value = 5642
$wsAll.Range("F46").Value = 4148
$wsAll.Range("F50").Value = 163
